$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "ns20232229@sav.edu.eg"
$ws.Range("B69").Value = "https://nesmasayed2004.github.io/mypage/"
